# Update the "From" greeting text for rule R10 (cell E8) from "Good Morning"
# to "GIT UPDATE", and leave the selection on that cell, matching the
# author's commit ("update file with jgit").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
[void]$ws.Range("E8").Select()
